# "Fruta / hortaliza, semanal"
# Insert a new weekly price-report row for Berenjena (Terminal Hortofrutícola
# Agro Chillán) at row 6, pushing the existing rows 6-18 down to rows 7-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 6..18 down to 7..19, leaving a blank row 6 (keeps
# formatting of the row it was inserted above, matching column D's date style).
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new weekly data point.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44600
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112001
$ws.Range("G6").Value = "Berenjena"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("N6").Value = "`$/caja 60 unidades"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 208
$ws.Range("Q6").Value = 60
$ws.Range("R6").Value = "Hortaliza"
